# Update the sheet with new daily data rows (245-247), matching the
# existing style/format used for the previous rows (date in col A uses the
# same style as the cell above it; B/C/D are plain numbers).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newRows = @(
    @{ Row = 245; A = 44319; B = 3; C = 40; D = 116.3941104580108 },
    @{ Row = 246; A = 44320; B = 0; C = 34; D = 98.9349938893092 },
    @{ Row = 247; A = 44321; B = 2; C = 32; D = 93.11528836640865 }
)

foreach ($r in $newRows) {
    $rowIndex = $r.Row

    # Column A carries the same date/border/bold style as the row above it;
    # copy that cell's formatting (and value, momentarily) then overwrite
    # the value with this row's date serial.
    $srcA = $ws.Cells.Item($rowIndex - 1, 1)
    $cellA = $ws.Cells.Item($rowIndex, 1)
    $srcA.Copy($cellA)
    $cellA.Value = $r.A

    $ws.Cells.Item($rowIndex, 2).Value = $r.B
    $ws.Cells.Item($rowIndex, 3).Value = $r.C
    $ws.Cells.Item($rowIndex, 4).Value = $r.D
}
